$wb = $excel.ActiveWorkbook

# ALC!row15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1340.2329
$ws.Range("I15").Value = 1340.2329
$ws.Range("K15").Value = 4020.6987
$ws.Range("M15").Value = -3851.6987

# ALC!row107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1010.4
$ws.Range("I107").Value = 1259.091
$ws.Range("J107").Value = 326.5
$ws.Range("K107").Value = 1259.091
$ws.Range("L107").Value = 326.5
$ws.Range("M107").Value = 660.9090000000001
$ws.Range("N107").Value = -4166.5

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 983.3333
$ws.Range("I137").Value = 842.6316
$ws.Range("J137").Value = 1518
$ws.Range("K137").Value = 2527.8948
$ws.Range("L137").Value = 4554
$ws.Range("M137").Value = 22.10519999999997
$ws.Range("N137").Value = -9654

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2994.638
$ws.Range("I138").Value = 2085.9546
$ws.Range("J138").Value = 3549.9443
$ws.Range("K138").Value = 6257.8638
$ws.Range("L138").Value = 10649.8329
$ws.Range("M138").Value = -1117.8638
$ws.Range("N138").Value = -20929.8329

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 3725.6924
$ws.Range("I141").Value = 3760.3809
$ws.Range("K141").Value = 11281.1427
$ws.Range("M141").Value = -6101.1427

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 428016.1
$ws.Range("I32").Value = 4293.418
$ws.Range("K32").Value = 4293.418
$ws.Range("M32").Value = -4006.418

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1301.9231
$ws.Range("I61").Value = 1251.3
$ws.Range("J61").Value = 1470.6666
$ws.Range("K61").Value = 1251.3
$ws.Range("L61").Value = 1470.6666
$ws.Range("M61").Value = -1039.3
$ws.Range("N61").Value = -1894.6666

# ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3761.45
$ws.Range("I88").Value = 2682.75
$ws.Range("J88").Value = 4480.5835
$ws.Range("K88").Value = 2682.75
$ws.Range("L88").Value = 4480.5835
$ws.Range("M88").Value = -2276.75
$ws.Range("N88").Value = -5292.5835

# ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 3761.45
$ws.Range("I91").Value = 2682.75
$ws.Range("J91").Value = 4480.5835
$ws.Range("K91").Value = 2682.75
$ws.Range("L91").Value = 4480.5835
$ws.Range("M91").Value = -1278.75
$ws.Range("N91").Value = -7288.5835

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 34599.098
$ws.Range("I122").Value = 2502.6
$ws.Range("J122").Value = 92956.37
$ws.Range("K122").Value = 7507.799999999999
$ws.Range("L122").Value = 278869.11
$ws.Range("M122").Value = -5057.799999999999
$ws.Range("N122").Value = -283769.11

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 15401388
$ws.Range("I132").Value = 18182936
$ws.Range("J132").Value = 102871.4
$ws.Range("K132").Value = 54548808
$ws.Range("L132").Value = 308614.2
$ws.Range("M132").Value = -54546278
$ws.Range("N132").Value = -313674.2

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1301.9231
$ws.Range("I136").Value = 1251.3
$ws.Range("J136").Value = 1470.6666
$ws.Range("K136").Value = 3753.9
$ws.Range("L136").Value = 4411.9998
$ws.Range("M136").Value = -1203.9
$ws.Range("N136").Value = -9511.9998

# BSM!row86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 35752308
$ws.Range("I86").Value = 58825812
$ws.Range("J86").Value = 93256.55
$ws.Range("K86").Value = 58825812
$ws.Range("L86").Value = 93256.55
$ws.Range("M86").Value = -58824689
$ws.Range("N86").Value = -95502.55

# BSM!row89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 35752308
$ws.Range("I89").Value = 58825812
$ws.Range("J89").Value = 93256.55
$ws.Range("K89").Value = 294129060
$ws.Range("L89").Value = 466282.75
$ws.Range("M89").Value = -294123444
$ws.Range("N89").Value = -477514.75

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5477.7896
$ws.Range("I134").Value = 2136.4783
$ws.Range("J134").Value = 10601.134
$ws.Range("K134").Value = 6409.4349
$ws.Range("L134").Value = 31803.402
$ws.Range("M134").Value = -3874.4349
$ws.Range("N134").Value = -36873.402

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7579790
$ws.Range("I31").Value = 12501905
$ws.Range("J31").Value = 7304.923
$ws.Range("K31").Value = 12501905
$ws.Range("L31").Value = 7304.923
$ws.Range("M31").Value = -12501610
$ws.Range("N31").Value = -7894.923

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 7579790
$ws.Range("I34").Value = 12501905
$ws.Range("J34").Value = 7304.923
$ws.Range("K34").Value = 12501905
$ws.Range("L34").Value = 7304.923
$ws.Range("M34").Value = -12501703
$ws.Range("N34").Value = -7708.923

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 447.66666
$ws.Range("I107").Value = 382.75
$ws.Range("J107").Value = 534.2222
$ws.Range("K107").Value = 382.75
$ws.Range("L107").Value = 534.2222
$ws.Range("M107").Value = 1537.25
$ws.Range("N107").Value = -4374.2222

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1531.1428
$ws.Range("I122").Value = 1579.6364
$ws.Range("J122").Value = 1477.8
$ws.Range("K122").Value = 4738.9092
$ws.Range("L122").Value = 4433.4
$ws.Range("M122").Value = -2288.9092
$ws.Range("N122").Value = -9333.4

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 501932.53
$ws.Range("I132").Value = 78143.234
$ws.Range("K132").Value = 234429.702
$ws.Range("M132").Value = -231899.702

# GSM!row136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 28166.666
$ws.Range("J136").Value = 28166.666
$ws.Range("L136").Value = 84499.99800000001
$ws.Range("N136").Value = -89599.99800000001

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3085.0625
$ws.Range("J7").Value = 3335.5386
$ws.Range("L7").Value = 3335.5386
$ws.Range("N7").Value = -3559.5386

# LTW!row22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1149.8572
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1149.8572
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1149.8572
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1739.8572

# LTW!row27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1149.8572
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1149.8572
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1149.8572
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1363.8572

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2056.639
$ws.Range("I40").Value = 1837.4166
$ws.Range("J40").Value = 2495.0833
$ws.Range("K40").Value = 1837.4166
$ws.Range("L40").Value = 2495.0833
$ws.Range("M40").Value = -1701.4166
$ws.Range("N40").Value = -2767.0833

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6947.1177
$ws.Range("I46").Value = 1011.2222
$ws.Range("J46").Value = 13625
$ws.Range("K46").Value = 1011.2222
$ws.Range("L46").Value = 13625
$ws.Range("M46").Value = -823.2222
$ws.Range("N46").Value = -14001

# LTW!row68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1990.4878
$ws.Range("I68").Value = 1925.1428
$ws.Range("J68").Value = 2131.2307
$ws.Range("K68").Value = 1925.1428
$ws.Range("L68").Value = 2131.2307
$ws.Range("M68").Value = -1176.1428
$ws.Range("N68").Value = -3629.2307

# LTW!row71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1990.4878
$ws.Range("I71").Value = 1925.1428
$ws.Range("J71").Value = 2131.2307
$ws.Range("K71").Value = 9625.714
$ws.Range("L71").Value = 10656.1535
$ws.Range("M71").Value = -5881.714
$ws.Range("N71").Value = -18144.1535

# LTW!row93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5614.143
$ws.Range("I93").Value = 7626.5
$ws.Range("J93").Value = 2931
$ws.Range("K93").Value = 7626.5
$ws.Range("L93").Value = 2931
$ws.Range("M93").Value = -6378.5
$ws.Range("N93").Value = -5427

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2279.1155
$ws.Range("I122").Value = 2069.25
$ws.Range("J122").Value = 2614.9
$ws.Range("K122").Value = 6207.75
$ws.Range("L122").Value = 7844.700000000001
$ws.Range("M122").Value = -3757.75
$ws.Range("N122").Value = -12744.7

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3085.0625
$ws.Range("J126").Value = 3335.5386
$ws.Range("L126").Value = 10006.6158
$ws.Range("N126").Value = -14946.6158

# WVR!row62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5432.579
$ws.Range("I62").Value = 5629.7144
$ws.Range("K62").Value = 5629.7144
$ws.Range("M62").Value = -5005.7144

# WVR!row65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5432.579
$ws.Range("I65").Value = 5629.7144
$ws.Range("K65").Value = 28148.572
$ws.Range("M65").Value = -25028.572

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3341.9048
$ws.Range("I81").Value = 2758
$ws.Range("J81").Value = 4801.6665
$ws.Range("K81").Value = 5516
$ws.Range("L81").Value = 9603.333000000001
$ws.Range("M81").Value = -4455
$ws.Range("N81").Value = -11725.333

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 3341.9048
$ws.Range("I84").Value = 2758
$ws.Range("J84").Value = 4801.6665
$ws.Range("K84").Value = 27580
$ws.Range("L84").Value = 48016.665
$ws.Range("M84").Value = -22276
$ws.Range("N84").Value = -58624.665

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 51197.75
$ws.Range("I136").Value = 53813.42
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 161440.26
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -158890.26
$ws.Range("N136").Value = -9600
